$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 153.8125
$ws.Range("I12").Value = 149.33333
$ws.Range("K12").Value = 149.33333
$ws.Range("M12").Value = 20.66667000000001
$ws.Range("H38").Value = 7105.2593
$ws.Range("I38").Value = 7928.5625
$ws.Range("K38").Value = 23785.6875
$ws.Range("M38").Value = -23413.6875
$ws.Range("H75").Value = 333372800
$ws.Range("J75").Value = 333372800
$ws.Range("L75").Value = 333372800
$ws.Range("N75").Value = -333374672
$ws.Range("H78").Value = 333372800
$ws.Range("J78").Value = 333372800
$ws.Range("L78").Value = 1000118400
$ws.Range("N78").Value = -1000127760
$ws.Range("H88").Value = 3429.9285
$ws.Range("J88").Value = 3861.2917
$ws.Range("L88").Value = 3861.2917
$ws.Range("N88").Value = -4673.2917
$ws.Range("H91").Value = 3429.9285
$ws.Range("J91").Value = 3861.2917
$ws.Range("L91").Value = 3861.2917
$ws.Range("N91").Value = -6669.2917
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6508.346
$ws.Range("I32").Value = 3941.5908
$ws.Range("K32").Value = 3941.5908
$ws.Range("M32").Value = -3654.5908
$ws.Range("H61").Value = 5412.0713
$ws.Range("I61").Value = 5063.0386
$ws.Range("K61").Value = 5063.0386
$ws.Range("M61").Value = -4851.0386
$ws.Range("H102").Value = 4625.926
$ws.Range("I102").Value = 1348.125
$ws.Range("J102").Value = 9393.637000000001
$ws.Range("K102").Value = 1348.125
$ws.Range("L102").Value = 9393.637000000001
$ws.Range("M102").Value = 273.875
$ws.Range("N102").Value = -12637.637
$ws.Range("H122").Value = 3836347.5
$ws.Range("I122").Value = 5849155
$ws.Range("K122").Value = 17547465
$ws.Range("M122").Value = -17545015
$ws.Range("H132").Value = 2506.3914
$ws.Range("I132").Value = 1649.8572
$ws.Range("K132").Value = 4949.571599999999
$ws.Range("M132").Value = -2419.571599999999
$ws.Range("H136").Value = 5412.0713
$ws.Range("I136").Value = 5063.0386
$ws.Range("K136").Value = 15189.1158
$ws.Range("M136").Value = -12639.1158
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 19999
$ws.Range("J46").Value = 19999
$ws.Range("L46").Value = 19999
$ws.Range("N46").Value = -20595
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4272.6
$ws.Range("I31").Value = 1829.8572
$ws.Range("K31").Value = 1829.8572
$ws.Range("M31").Value = -1534.8572
$ws.Range("H34").Value = 4272.6
$ws.Range("I34").Value = 1829.8572
$ws.Range("K34").Value = 1829.8572
$ws.Range("M34").Value = -1627.8572
$ws.Range("H50").Value = 54997.5
$ws.Range("J50").Value = 54997.5
$ws.Range("L50").Value = 54997.5
$ws.Range("N50").Value = -56247.5
$ws.Range("H51").Value = 39999.5
$ws.Range("J51").Value = 39999.5
$ws.Range("L51").Value = 39999.5
$ws.Range("N51").Value = -41471.5
$ws.Range("H61").Value = 39999.5
$ws.Range("J61").Value = 39999.5
$ws.Range("L61").Value = 39999.5
$ws.Range("N61").Value = -40695.5
$ws.Range("H107").Value = 1370.1
$ws.Range("I107").Value = 1423.5555
$ws.Range("J107").Value = 889
$ws.Range("K107").Value = 1423.5555
$ws.Range("L107").Value = 889
$ws.Range("M107").Value = 496.4445000000001
$ws.Range("H132").Value = 2483.7334
$ws.Range("I132").Value = 2089.7144
$ws.Range("K132").Value = 6269.1432
$ws.Range("M132").Value = -3739.1432
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 275325.94
$ws.Range("I11").Value = 510818.5
$ws.Range("J11").Value = 584.5833
$ws.Range("K11").Value = 1532455.5
$ws.Range("L11").Value = 1753.7499
$ws.Range("M11").Value = -1532315.5
$ws.Range("N11").Value = -2033.7499
$ws.Range("H39").Value = 1984.6471
$ws.Range("J39").Value = 1701.1428
$ws.Range("L39").Value = 5103.428400000001
$ws.Range("N39").Value = -5691.428400000001
$ws.Range("H101").Value = 9718.833000000001
$ws.Range("J101").Value = 9718.833000000001
$ws.Range("L101").Value = 29156.499
$ws.Range("N101").Value = -34024.499
$ws.Range("H114").Value = 366.33334
$ws.Range("I114").Value = 366.33334
$ws.Range("K114").Value = 1099.00002
$ws.Range("M114").Value = 2154.99998
$ws.Range("H121").Value = 1653.2
$ws.Range("J121").Value = 2347.1052
$ws.Range("L121").Value = 7041.3156
$ws.Range("N121").Value = -9661.3156
$ws.Range("H137").Value = 13992.4
$ws.Range("J137").Value = 16943.666
$ws.Range("L137").Value = 50830.99800000001
$ws.Range("N137").Value = -61030.99800000001
$ws.Range("H138").Value = 33336120
$ws.Range("I138").Value = 55556924
$ws.Range("J138").Value = 4916
$ws.Range("K138").Value = 166670772
$ws.Range("L138").Value = 14748
$ws.Range("M138").Value = -166665632
$ws.Range("N138").Value = -25028
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 8162222
$ws.Range("J11").Value = 33744.832
$ws.Range("L11").Value = 33744.832
$ws.Range("N11").Value = -34022.832
$ws.Range("H102").Value = 1461.9
$ws.Range("I102").Value = 577.375
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 577.375
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = 1044.625
$ws.Range("N102").Value = -8244
$ws.Range("H113").Value = 27785854
$ws.Range("I113").Value = 52639680
$ws.Range("K113").Value = 52639680
$ws.Range("M113").Value = -52637510
$ws.Range("H122").Value = 2622.9473
$ws.Range("I122").Value = 1427.25
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 4281.75
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -1831.75
$ws.Range("N122").Value = -31900
$ws.Range("H126").Value = 4225
$ws.Range("I126").Value = 2837.5
$ws.Range("K126").Value = 8512.5
$ws.Range("M126").Value = -6042.5
$ws.Range("H132").Value = 7405.909
$ws.Range("I132").Value = 2727.8333
$ws.Range("J132").Value = 13019.6
$ws.Range("K132").Value = 8183.499899999999
$ws.Range("L132").Value = 39058.8
$ws.Range("M132").Value = -5653.499899999999
$ws.Range("N132").Value = -44118.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4199.5
$ws.Range("I7").Value = 4085.1428
$ws.Range("K7").Value = 4085.1428
$ws.Range("M7").Value = -3973.1428
$ws.Range("H40").Value = 3362.818
$ws.Range("I40").Value = 797.4286
$ws.Range("K40").Value = 797.4286
$ws.Range("M40").Value = -661.4286
$ws.Range("H126").Value = 4199.5
$ws.Range("I126").Value = 4085.1428
$ws.Range("K126").Value = 12255.4284
$ws.Range("M126").Value = -9785.428400000001
$ws.Range("H132").Value = 4931.077
$ws.Range("I132").Value = 4316.5
$ws.Range("K132").Value = 12949.5
$ws.Range("M132").Value = -10419.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 95038.2
$ws.Range("J93").Value = 95038.2
$ws.Range("L93").Value = 95038.2
$ws.Range("N93").Value = -100030.2
$ws.Range("H132").Value = 3314.375
$ws.Range("I132").Value = 3104.3
$ws.Range("K132").Value = 9312.900000000001
$ws.Range("M132").Value = -6782.900000000001
$ws.Range("H136").Value = 6660.8667
$ws.Range("J136").Value = 5991.5
$ws.Range("L136").Value = 17974.5
$ws.Range("N136").Value = -23074.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N107").Value = -4729
